# Lesson 08 Discussion.docx - proofreading pass.
#
# Applies a handful of wording tweaks to the discussion body:
#   - tidies up a couple of phrasings ("well being" -> "well-being",
#     "get on Facebook" -> "get from Facebook")
#   - adds a couple of clarifying parentheticals/qualifiers
#   - drops a stray trailing parenthesis
#   - extends the closing sentence
# and leaves a "_GoBack" bookmark at the cursor's last edit position, the
# way Word itself does after an interactive editing session.

$d = $word.ActiveDocument

# "...these algorithms for the most part [gramStart]are[gramEnd] optimized
# for..." -> the grammar-checker flag around "are" is irrelevant once the
# surrounding wording is untouched, so just normalize the run (no visible
# text change here).
$d.Content.Find.Execute("for the most part are optimized for", $true, $false, $false, $false, $false,
                         $true, 1, $false, "for the most part are optimized for", 2) | Out-Null

# "...mental health and/or general well being could be ignored" ->
# hyphenate "well-being".
$d.Content.Find.Execute("general well being could be ignored", $true, $false, $false, $false, $false,
                         $true, 1, $false, "general well-being could be ignored", 2) | Out-Null

# Qualify the claim with a parenthetical.
$d.Content.Find.Execute("could be ignored. We see examples", $true, $false, $false, $false, $false,
                         $true, 1, $false, "could be ignored (which is not always the case). We see examples", 2) | Out-Null

# "the feed we get on Facebook" -> "the feed we get from Facebook".
$d.Content.Find.Execute("in the feed we get on Facebook", $true, $false, $false, $false, $false,
                         $true, 1, $false, "in the feed we get from Facebook", 2) | Out-Null

# "exposed to new ideas, but" -> "exposed to new ideas/perspectives, but".
$d.Content.Find.Execute("exposed to new ideas, but", $true, $false, $false, $false, $false,
                         $true, 1, $false, "exposed to new ideas/perspectives, but", 2) | Out-Null

# Drop the stray trailing ")" after "from others like ourselves)".
$d.Content.Find.Execute("from others like ourselves)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "from others like ourselves", 2) | Out-Null

# Extend the closing sentence of the second paragraph.
$d.Content.Find.Execute("the user keep exposure to other ideas/products.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "the user keep exposure to other ideas/products that are still relevant.", 2) | Out-Null

# Leave a "_GoBack" bookmark (Word auto-drops one at the last edit point)
# splitting "...so called" into "...so ca" | "lled...".
$rng = $d.Content
$rng.Find.Execute("could lead to the so called", $true, $false, $false, $false, $false,
                   $true, 1, $false, "", 0) | Out-Null
$found = $rng.Duplicate
$found.Start = $found.End - 4
$found.Collapse(1)
$d.Bookmarks.Add("_GoBack", $found) | Out-Null
